# SLE 2024-12 Closed Bugs update to include 2935004
# 2935004: Fix ServiceRegistry not being able to communicate with JupyterHub
#          when a CNI with NetworkPolicy support is installed

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("closed bugs in last iteration")

# Correct the wording of the existing "loading mask" bug title (row 5 / ID 2929769)
$ws.Range("B5").Value = "Loading mask glitch when deleting more packages in Packages grid from Feed details"

# Rename the Grafana dependency bug title (row 12 / ID 2901954)
$ws.Range("B12").Value = "Security vulnerabilities in ni-grafana"

# Append the new closed bug as row 17
$ws.Range("A17").Value = 2935004
$ws.Range("B17").Value = "Fix ServiceRegistry not being able to communicate with JupyterHub when a CNI with NetworkPolicy support is installed"
$ws.Range("C17").Value = "Closed"

# Match the formatting of the preceding data rows
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection to A17, where the new row was entered
$ws.Range("A17").Select()
